$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (changed) date column C was bumped by one day
# (45181 -> 45182, i.e. 2023-09-12 -> 2023-09-13) for every data row
# (rows 2 through 407).
$ws.Range("C2:C407").Value = 45182
